# Actualización automática desde WSL
# Refresh the weather-station data sheet: correct the sub-millisecond
# rounding of the last existing timestamp (row 9) and append the new
# reading captured for 09:00 (row 10).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 9: timestamp got re-serialised with millisecond precision ---
$ws.Range("A9").Value = 45878.33351998842

# --- Row 10: new hourly reading ---
$ws.Range("A10").Value = 45878.37517072909
$ws.Range("A10").NumberFormat = $ws.Range("A9").NumberFormat

$ws.Range("B10").Value = 2025
$ws.Range("C10").Value = 37
$ws.Range("D10").Value = 13.89
$ws.Range("E10").Value = 92.75
$ws.Range("F10").Value = 161.72
$ws.Range("G10").Value = 4.31
$ws.Range("H10").Value = "ESE"
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = "09:00:14"
